$wb = $excel.ActiveWorkbook

# Rename the existing sheet from "Tabelle1" to "Tests"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Tests"

# Add a new worksheet named "Gerätespecs" after the first sheet
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Gerätespecs"

# Re-activate the Tests sheet and scroll so column B is the leftmost visible column
$ws1.Activate()
$excel.ActiveWindow.ScrollColumn = 2

